# Update "Forecast Comparison" sheet with a new Week_Start_Date column
# and corrected forecast numbers, then update the "Summary" sheet totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")
$summary = $wb.Worksheets.Item("Summary")

# ------------------------------------------------------------------
# 1. Insert a new column B ("Week_Start_Date") - this shifts the
#    existing ASIN / MyForecast / Amazon*/ Product Title /
#    is_holiday_week columns one place to the right (B->C, C->D, ...).
# ------------------------------------------------------------------
$ws.Columns.Item(2).Insert()

# ------------------------------------------------------------------
# 2. Header for the new column
# ------------------------------------------------------------------
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# ------------------------------------------------------------------
# 3. Per-row data: Week label (drop leading zero), Week_Start_Date,
#    and corrected MyForecast values.
# ------------------------------------------------------------------
$weeks = @(1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16)
$weekStartDates = @(
    "2025-01-05", "2025-01-12", "2025-01-19", "2025-01-26",
    "2025-02-02", "2025-02-09", "2025-02-16", "2025-02-23",
    "2025-03-02", "2025-03-09", "2025-03-16", "2025-03-23",
    "2025-03-30", "2025-04-06", "2025-04-13", "2025-04-20"
)
$myForecast = @(67,69,70,71,71,74,77,78,77,74,72,70,68,69,72,72)

# make sure the new date column is plain text so Excel does not
# auto-convert the "yyyy-mm-dd" strings into date serial numbers
$ws.Range("B2:B17").NumberFormat = "@"

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2

    # Week label without the leading zero ("W01" -> "W1", ... "W9")
    $ws.Cells.Item($row, 1).Value = "W" + $weeks[$i]

    # New Week_Start_Date column
    $ws.Cells.Item($row, 2).Value = $weekStartDates[$i]

    # Corrected MyForecast value (column D after the insert)
    $ws.Cells.Item($row, 4).Value = $myForecast[$i]

    # is_holiday_week (column J after the insert) becomes a boolean
    $ws.Cells.Item($row, 10).Value = $false
}

# ------------------------------------------------------------------
# 4. Update the Summary sheet totals that depend on MyForecast.
#    Keep them as plain text (matching the rest of the "Value"
#    column on this sheet) rather than letting Excel coerce them
#    into numbers.
# ------------------------------------------------------------------
$summary.Range("B9").NumberFormat  = "@"
$summary.Range("B10").NumberFormat = "@"
$summary.Range("B11").NumberFormat = "@"
$summary.Range("B12").NumberFormat = "@"
$summary.Range("B14").NumberFormat = "@"

$summary.Range("B9").Value  = "1151"   # Total Forecast (16 Weeks)
$summary.Range("B10").Value = "577"    # Total Forecast (8 Weeks)
$summary.Range("B11").Value = "277"    # Total Forecast (4 Weeks)
$summary.Range("B12").Value = "78"     # Max Forecast
$summary.Range("B14").Value = "67"     # Min Forecast
